# The deck's theme (ppt/theme/theme1.xml, used by the one slide master /
# all slides) currently carries the "Integral" color scheme. The commit
# swaps it for the stock "Office Theme" color scheme (the 12 standard
# Office theme colors), while the secondary theme part used by the notes
# master ends up holding the colors the primary theme used to have.
#
# The PowerPoint object model exposes the 12 scheme colors (dk1, lt1,
# dk2, lt2, accent1-6, hlink, folHlink) through ThemeColorScheme, which
# is reachable from a Slide (and therefore resolves against the deck's
# one-and-only slide master / theme part). We drive that to repaint the
# theme's color scheme to the standard Office values.

$p = $ppt.ActivePresentation

# Standard "Office Theme" color scheme, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
# accent6, hlink, folHlink
$officeColors = @(
    @(0x00, 0x00, 0x00),  # dk1
    @(0xFF, 0xFF, 0xFF),  # lt1
    @(0x44, 0x54, 0x6A),  # dk2
    @(0xE7, 0xE6, 0xE6),  # lt2
    @(0x5B, 0x9B, 0xD5),  # accent1
    @(0xED, 0x7D, 0x31),  # accent2
    @(0xA5, 0xA5, 0xA5),  # accent3
    @(0xFF, 0xC0, 0x00),  # accent4
    @(0x44, 0x72, 0xC4),  # accent5
    @(0x70, 0xAD, 0x47),  # accent6
    @(0x05, 0x63, 0xC1),  # hlink
    @(0x95, 0x4F, 0x72)   # folHlink
)

$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $rgbTriplet = $officeColors[$i]
    $r = $rgbTriplet[0]
    $g = $rgbTriplet[1]
    $b = $rgbTriplet[2]
    $oleColor = $r + ($g * 256) + ($b * 65536)
    $tcs.Item($i + 1).RGB = $oleColor
}
